$d = $word.ActiveDocument

$pairs = @(
    @("55×26=", "36×76="),
    @("56×77=", "93×92="),
    @("38×98=", "93×31="),
    @("29×26=", "20×23="),
    @("39×56=", "28×53="),
    @("29×17=", "23×98="),
    @("19×42=", "30×49="),
    @("15×50=", "66×51="),
    @("57×90=", "79×82="),
    @("68×66=", "54×76="),
    @("58×75=", "48×93="),
    @("11×77=", "54×39="),
    @("84×99=", "65×44="),
    @("75×53=", "26×18="),
    @("79×13=", "24×35="),
    @("88×89=", "45×35="),
    @("65×98=", "12×86="),
    @("94×15=", "93×28="),
    @("38×81=", "18×77="),
    @("39×87=", "51×88="),
    @("34×41=", "77×85="),
    @("80×97=", "33×78="),
    @("86×20=", "79×41="),
    @("18×33=", "30×69="),
    @("18×84=", "89×56=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
